# Apply the "constructive method and report finished" update:
#  - Resumen sheet: update best zone (Z2 -> Z1) and its max value
#  - Solucion sheet: replace the Pedido/Salida assignment table with the new solution
#  - Metricas sheet: update the Z1/Z2 time values

$wb = $excel.ActiveWorkbook

# ---- Resumen sheet ----
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("B2").Value = "Z1"
$wsResumen.Range("C2").Value = 514.8879554546438

# ---- Solucion sheet ----
$wsSolucion = $wb.Worksheets.Item("Solucion")

$solucion = @(
    @("Pedido_36", "S001"),
    @("Pedido_14", "S021"),
    @("Pedido_39", "S011"),
    @("Pedido_5", "S031"),
    @("Pedido_23", "S002"),
    @("Pedido_4", "S012"),
    @("Pedido_6", "S022"),
    @("Pedido_32", "S003"),
    @("Pedido_16", "S032"),
    @("Pedido_26", "S013"),
    @("Pedido_27", "S004"),
    @("Pedido_35", "S023"),
    @("Pedido_7", "S014"),
    @("Pedido_11", "S033"),
    @("Pedido_21", "S005"),
    @("Pedido_30", "S024"),
    @("Pedido_33", "S015"),
    @("Pedido_40", "S034"),
    @("Pedido_37", "S025"),
    @("Pedido_25", "S006"),
    @("Pedido_17", "S035"),
    @("Pedido_15", "S016"),
    @("Pedido_31", "S026"),
    @("Pedido_2", "S036"),
    @("Pedido_18", "S007"),
    @("Pedido_3", "S027"),
    @("Pedido_20", "S017"),
    @("Pedido_22", "S008"),
    @("Pedido_1", "S037"),
    @("Pedido_12", "S028"),
    @("Pedido_13", "S018"),
    @("Pedido_28", "S038"),
    @("Pedido_34", "S009"),
    @("Pedido_10", "S019"),
    @("Pedido_29", "S029"),
    @("Pedido_24", "S039"),
    @("Pedido_38", "S010"),
    @("Pedido_8", "S030"),
    @("Pedido_19", "S040"),
    @("Pedido_9", "S020")
)

for ($i = 0; $i -lt $solucion.Length; $i++) {
    $row = $i + 2
    $wsSolucion.Cells.Item($row, 1).Value = $solucion[$i][0]
    $wsSolucion.Cells.Item($row, 2).Value = $solucion[$i][1]
}

# ---- Metricas sheet ----
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Range("B2").Value = 514.8879554546438
$wsMetricas.Range("B3").Value = 513.1332630554655
